$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the "First Practice: https://.../474eb226db" paragraph that
# precedes the empty paragraph we need to fill in. We anchor on the
# unique URL fragment rather than a hard-coded paragraph index so the
# script keeps working even if unrelated paragraphs shift around.
# ------------------------------------------------------------------
$firstPracticeIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*474eb226db*") {
        $firstPracticeIndex = $i
        break
    }
}

$secondPracticeIndex = $firstPracticeIndex + 1
$tabParaIndex = $firstPracticeIndex + 3
$motivationIndex = $firstPracticeIndex + 5
$jobPostIndex = $firstPracticeIndex + 6

# ------------------------------------------------------------------
# 1. Fill in the empty paragraph right after "First Practice: ..."
#    with "Second Practie: " + hyperlink + trailing space.
# ------------------------------------------------------------------
$secondPracticePara = $d.Paragraphs.Item($secondPracticeIndex)
$secondPracticeRange = $secondPracticePara.Range
# Insert the text with a one-character marker at the end; that marker
# will be swapped out for the hyperlink in the next step.
$secondPracticeRange.InsertBefore("Second Practie: X")

$secondPracticePara = $d.Paragraphs.Item($secondPracticeIndex)
$secondPracticeRange = $secondPracticePara.Range
$markerRange = $d.Range($secondPracticeRange.End - 2, $secondPracticeRange.End - 1)
$d.Hyperlinks.Add($markerRange, "https://app.biginterview.com/s/6790893c9e", "", "", "https://app.biginterview.com/s/6790893c9e")

$secondPracticePara = $d.Paragraphs.Item($secondPracticeIndex)
$secondPracticeRange = $secondPracticePara.Range
$secondPracticeRange.InsertAfter(" ")

# ------------------------------------------------------------------
# 2. Remove the lone <w:tab/> run from the tab-stop paragraph, leaving
#    its paragraph properties (the tab stop definition) untouched.
# ------------------------------------------------------------------
$tabPara = $d.Paragraphs.Item($tabParaIndex)
$tabRange = $tabPara.Range
$tabCharRange = $d.Range($tabRange.Start, $tabRange.Start + 1)
$tabCharRange.Delete()

# ------------------------------------------------------------------
# 3. Split "... help out my parents" into "... help out my " + "family"
# ------------------------------------------------------------------
$motivationPara = $d.Paragraphs.Item($motivationIndex)
$searchRange = $d.Range($motivationPara.Range.Start, $motivationPara.Range.End)
$findResult = $searchRange.Find.Execute("parents", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$parentsStart = $searchRange.Start
$parentsEnd = $searchRange.End
$parentsRange = $d.Range($parentsStart, $parentsEnd)
$parentsRange.Delete()
$insertPoint = $d.Range($parentsStart, $parentsStart)
$insertPoint.InsertAfter("family")
$familyRunRange = $d.Range($parentsStart, $parentsStart + 6)
# Toggle a format on/off so this new text commits as its own run while
# ending up with formatting identical to its neighbouring run.
$familyRunRange.Bold = 1
$familyRunRange.Bold = 0

# ------------------------------------------------------------------
# 4. Replace the whole "I was looking for..." sentence with
#    "I was excited when I found this position on Linked" + "In"
# ------------------------------------------------------------------
$jobPostPara = $d.Paragraphs.Item($jobPostIndex)
$jobPostFull = $jobPostPara.Range
$jobPostTextRange = $d.Range($jobPostFull.Start, $jobPostFull.End - 1)
$jobPostTextRange.Text = "I was excited when I found this position on Linked"

$jobPostPara = $d.Paragraphs.Item($jobPostIndex)
$jobPostRange = $jobPostPara.Range
$inInsertPoint = $d.Range($jobPostRange.End - 1, $jobPostRange.End - 1)
$inInsertPoint.InsertAfter("In")

$jobPostPara = $d.Paragraphs.Item($jobPostIndex)
$jobPostRange = $jobPostPara.Range
$inRunRange = $d.Range($jobPostRange.End - 3, $jobPostRange.End - 1)
$inRunRange.Bold = 1
$inRunRange.Bold = 0
